# Insert a new weekly price record as row 125, pushing existing rows 125-169 down to 126-170.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(125).Insert()

$row = 125
$ws.Cells.Item($row, 1).Value = 2
$ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44917
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100109
$ws.Cells.Item($row, 8).Value = "Uva"
$ws.Cells.Item($row, 9).Value = 100109001
$ws.Cells.Item($row, 10).Value = "Uva"
$ws.Cells.Item($row, 11).Value = "Flame Seedless"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 700
$ws.Cells.Item($row, 14).Value = 9500
$ws.Cells.Item($row, 15).Value = 10000
$ws.Cells.Item($row, 16).Value = 9750
$ws.Cells.Item($row, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 19).Value = 975
$ws.Cells.Item($row, 20).Value = 10
